$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "51.227.09"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "  -1.78%  "
$ws.Cells.Item(2, 5).ClearFormats()
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.920.45"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "  -2.41%  "
$ws.Cells.Item(3, 5).ClearFormats()
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.998"
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "  -0.31%  "
$ws.Cells.Item(4, 5).ClearFormats()
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "373.44"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "  +5.37%  "
$ws.Cells.Item(5, 5).ClearFormats()
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "102.43"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "  -5.65%  "
$ws.Cells.Item(6, 5).ClearFormats()
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.543"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "  -3.70%  "
$ws.Cells.Item(7, 5).ClearFormats()
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "  -0.13%  "
$ws.Cells.Item(8, 5).ClearFormats()
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.588"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "  -5.02%  "
$ws.Cells.Item(9, 5).ClearFormats()
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "37.02"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "  -3.39%  "
$ws.Cells.Item(10, 5).ClearFormats()
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "  +0.37%  "
$ws.Cells.Item(11, 5).ClearFormats()
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0836"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "  -2.57%  "
$ws.Cells.Item(12, 5).ClearFormats()
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "18.34"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "  -5.15%  "
$ws.Cells.Item(13, 5).ClearFormats()
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "3.372.84"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "  -2.38%  "
$ws.Cells.Item(14, 5).ClearFormats()
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.44"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "  -3.47%  "
$ws.Cells.Item(15, 5).ClearFormats()
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "2.910.49"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "  -2.14%  "
$ws.Cells.Item(16, 5).ClearFormats()
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.929"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "  -7.94%  "
$ws.Cells.Item(17, 5).ClearFormats()
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "51.119.74"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "  -2.09%  "
$ws.Cells.Item(18, 5).ClearFormats()
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "3.27"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "  -6.44%  "
$ws.Cells.Item(19, 5).ClearFormats()
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.25"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "  -3.71%  "
$ws.Cells.Item(20, 5).ClearFormats()
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "12.92"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "  -5.37%  "
$ws.Cells.Item(21, 5).ClearFormats()
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.0₃0945"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "  -2.97%  "
$ws.Cells.Item(22, 5).ClearFormats()
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "68.31"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "  -1.76%  "
$ws.Cells.Item(23, 5).ClearFormats()
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "260.00"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "  -1.43%  "
$ws.Cells.Item(24, 5).ClearFormats()
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.70"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "  -1.26%  "
$ws.Cells.Item(25, 5).ClearFormats()
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.170"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "  -5.13%  "
$ws.Cells.Item(26, 5).ClearFormats()
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "  +0.01%  "
$ws.Cells.Item(27, 5).ClearFormats()
$ws.Cells.Item(28, 2).NumberFormat = "@"
$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 2).ClearFormats()
$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 3).ClearFormats()
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "25.73"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "  -4.48%  "
$ws.Cells.Item(28, 5).ClearFormats()
$ws.Cells.Item(29, 2).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = "Filecoin"
$ws.Cells.Item(29, 2).ClearFormats()
$ws.Cells.Item(29, 3).NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(29, 3).ClearFormats()
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "7.14"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "  -6.57%  "
$ws.Cells.Item(29, 5).ClearFormats()
$ws.Cells.Item(30, 2).NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = "RenderToken"
$ws.Cells.Item(30, 2).ClearFormats()
$ws.Cells.Item(30, 3).NumberFormat = "@"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(30, 3).ClearFormats()
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "6.62"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "  +7.60%  "
$ws.Cells.Item(30, 5).ClearFormats()
$ws.Cells.Item(31, 2).NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 2).ClearFormats()
$ws.Cells.Item(31, 3).NumberFormat = "@"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(31, 3).ClearFormats()
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.102"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "  -6.90%  "
$ws.Cells.Item(31, 5).ClearFormats()
$ws.Cells.Item(32, 2).NumberFormat = "@"
$ws.Cells.Item(32, 2).Value = "Cosmos"
$ws.Cells.Item(32, 2).ClearFormats()
$ws.Cells.Item(32, 3).NumberFormat = "@"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(32, 3).ClearFormats()
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "9.88"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "  -4.26%  "
$ws.Cells.Item(32, 5).ClearFormats()
$ws.Cells.Item(33, 2).NumberFormat = "@"
$ws.Cells.Item(33, 2).Value = "Toncoin"
$ws.Cells.Item(33, 2).ClearFormats()
$ws.Cells.Item(33, 3).NumberFormat = "@"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(33, 3).ClearFormats()
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.11"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "  -3.66%  "
$ws.Cells.Item(33, 5).ClearFormats()
$ws.Cells.Item(34, 2).NumberFormat = "@"
$ws.Cells.Item(34, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(34, 2).ClearFormats()
$ws.Cells.Item(34, 3).NumberFormat = "@"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(34, 3).ClearFormats()
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "34.47"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = "  -5.04%  "
$ws.Cells.Item(34, 5).ClearFormats()
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "51.32"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "  +1.00%  "
$ws.Cells.Item(35, 5).ClearFormats()
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(36, 2).ClearFormats()
$ws.Cells.Item(36, 3).NumberFormat = "@"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(36, 3).ClearFormats()
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.00"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "  +0.26%  "
$ws.Cells.Item(36, 5).ClearFormats()
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = "VeChain"
$ws.Cells.Item(37, 2).ClearFormats()
$ws.Cells.Item(37, 3).NumberFormat = "@"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(37, 3).ClearFormats()
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0424"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "  -4.11%  "
$ws.Cells.Item(37, 5).ClearFormats()
$ws.Cells.Item(38, 2).NumberFormat = "@"
$ws.Cells.Item(38, 2).Value = "LidoDAOToken"
$ws.Cells.Item(38, 2).ClearFormats()
$ws.Cells.Item(38, 3).NumberFormat = "@"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(38, 3).ClearFormats()
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.99"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "  -6.61%  "
$ws.Cells.Item(38, 5).ClearFormats()
$ws.Cells.Item(39, 2).NumberFormat = "@"
$ws.Cells.Item(39, 2).Value = "Celestia"
$ws.Cells.Item(39, 2).ClearFormats()
$ws.Cells.Item(39, 3).NumberFormat = "@"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(39, 3).ClearFormats()
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "17.08"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "  -4.52%  "
$ws.Cells.Item(39, 5).ClearFormats()
$ws.Cells.Item(40, 2).NumberFormat = "@"
$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 2).ClearFormats()
$ws.Cells.Item(40, 3).NumberFormat = "@"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(40, 3).ClearFormats()
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.58"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "  -4.71%  "
$ws.Cells.Item(40, 5).ClearFormats()
$ws.Cells.Item(41, 2).NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = "ARBITRUM"
$ws.Cells.Item(41, 2).ClearFormats()
$ws.Cells.Item(41, 3).NumberFormat = "@"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(41, 3).ClearFormats()
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.84"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "  -6.75%  "
$ws.Cells.Item(41, 5).ClearFormats()
$ws.Cells.Item(42, 2).NumberFormat = "@"
$ws.Cells.Item(42, 2).Value = "Stellar"
$ws.Cells.Item(42, 2).ClearFormats()
$ws.Cells.Item(42, 3).NumberFormat = "@"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(42, 3).ClearFormats()
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.113"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "  -3.77%  "
$ws.Cells.Item(42, 5).ClearFormats()
$ws.Cells.Item(43, 2).NumberFormat = "@"
$ws.Cells.Item(43, 2).Value = "EnergySwap"
$ws.Cells.Item(43, 2).ClearFormats()
$ws.Cells.Item(43, 3).NumberFormat = "@"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(43, 3).ClearFormats()
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "22.11"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "  -2.93%  "
$ws.Cells.Item(43, 5).ClearFormats()
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "119.70"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "  -2.00%  "
$ws.Cells.Item(44, 5).ClearFormats()
$ws.Cells.Item(45, 2).NumberFormat = "@"
$ws.Cells.Item(45, 2).Value = "WEMIXToken"
$ws.Cells.Item(45, 2).ClearFormats()
$ws.Cells.Item(45, 3).NumberFormat = "@"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(45, 3).ClearFormats()
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.09"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "  -2.40%  "
$ws.Cells.Item(45, 5).ClearFormats()
$ws.Cells.Item(46, 2).NumberFormat = "@"
$ws.Cells.Item(46, 2).Value = "Maker"
$ws.Cells.Item(46, 2).ClearFormats()
$ws.Cells.Item(46, 3).NumberFormat = "@"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(46, 3).ClearFormats()
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.019.55"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "  -5.00%  "
$ws.Cells.Item(46, 5).ClearFormats()
$ws.Cells.Item(47, 2).NumberFormat = "@"
$ws.Cells.Item(47, 2).Value = "ApeXProtocol"
$ws.Cells.Item(47, 2).ClearFormats()
$ws.Cells.Item(47, 3).NumberFormat = "@"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(47, 3).ClearFormats()
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.31"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "  -3.07%  "
$ws.Cells.Item(47, 5).ClearFormats()
$ws.Cells.Item(48, 2).NumberFormat = "@"
$ws.Cells.Item(48, 2).Value = "NEARProtocol"
$ws.Cells.Item(48, 2).ClearFormats()
$ws.Cells.Item(48, 3).NumberFormat = "@"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(48, 3).ClearFormats()
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.16"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "  -6.10%  "
$ws.Cells.Item(48, 5).ClearFormats()
$ws.Cells.Item(49, 2).NumberFormat = "@"
$ws.Cells.Item(49, 2).Value = "RocketPoolETH"
$ws.Cells.Item(49, 2).ClearFormats()
$ws.Cells.Item(49, 3).NumberFormat = "@"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(49, 3).ClearFormats()
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.205.27"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "  -2.21%  "
$ws.Cells.Item(49, 5).ClearFormats()
$ws.Cells.Item(50, 2).NumberFormat = "@"
$ws.Cells.Item(50, 2).Value = "TheGraph"
$ws.Cells.Item(50, 2).ClearFormats()
$ws.Cells.Item(50, 3).NumberFormat = "@"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(50, 3).ClearFormats()
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.241"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "  +0.34%  "
$ws.Cells.Item(50, 5).ClearFormats()
$ws.Cells.Item(51, 2).NumberFormat = "@"
$ws.Cells.Item(51, 2).Value = "BEAM"
$ws.Cells.Item(51, 2).ClearFormats()
$ws.Cells.Item(51, 3).NumberFormat = "@"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Cells.Item(51, 3).ClearFormats()
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0311"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "  -7.22%  "
$ws.Cells.Item(51, 5).ClearFormats()
